$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 153789
$ws.Range("C4").Value = 145288
$ws.Range("C5").Value = 8501
$ws.Range("C7").Value = 5.53
$ws.Range("C8").Value = 63.81
